# Auto-generated script to apply scheduled market-data refresh values
# to the leveling-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 849.8570999999999
$ws.Cells.Item(58, 9).Value = 786.4
$ws.Cells.Item(58, 10).Value = 1008.5
$ws.Cells.Item(58, 11).Value = 2359.2
$ws.Cells.Item(58, 12).Value = 3025.5
$ws.Cells.Item(58, 13).Value = -2209.2
$ws.Cells.Item(58, 14).Value = -3325.5
$ws.Cells.Item(69, 8).Value = 8754.315000000001
$ws.Cells.Item(69, 9).Value = 5473
$ws.Cells.Item(69, 10).Value = 9629.333000000001
$ws.Cells.Item(69, 11).Value = 16419
$ws.Cells.Item(69, 12).Value = 28887.999
$ws.Cells.Item(69, 13).Value = -15545
$ws.Cells.Item(69, 14).Value = -30635.999
$ws.Cells.Item(72, 8).Value = 8754.315000000001
$ws.Cells.Item(72, 9).Value = 5473
$ws.Cells.Item(72, 10).Value = 9629.333000000001
$ws.Cells.Item(72, 11).Value = 49257
$ws.Cells.Item(72, 12).Value = 86663.997
$ws.Cells.Item(72, 13).Value = -44889
$ws.Cells.Item(72, 14).Value = -95399.997
$ws.Cells.Item(76, 8).Value = 3008.5454
$ws.Cells.Item(76, 10).Value = 995
$ws.Cells.Item(76, 12).Value = 995
$ws.Cells.Item(76, 14).Value = -1625
$ws.Cells.Item(79, 8).Value = 3008.5454
$ws.Cells.Item(79, 10).Value = 995
$ws.Cells.Item(79, 12).Value = 995
$ws.Cells.Item(79, 14).Value = -3179
$ws.Cells.Item(80, 8).Value = 4311363
$ws.Cells.Item(80, 9).Value = 6579907.5
$ws.Cells.Item(80, 10).Value = 1127.7
$ws.Cells.Item(80, 11).Value = 19739722.5
$ws.Cells.Item(80, 12).Value = 3383.1
$ws.Cells.Item(80, 13).Value = -19738724.5
$ws.Cells.Item(80, 14).Value = -5379.1
$ws.Cells.Item(83, 8).Value = 4311363
$ws.Cells.Item(83, 9).Value = 6579907.5
$ws.Cells.Item(83, 10).Value = 1127.7
$ws.Cells.Item(83, 11).Value = 59219167.5
$ws.Cells.Item(83, 12).Value = 10149.3
$ws.Cells.Item(83, 13).Value = -59214175.5
$ws.Cells.Item(83, 14).Value = -20133.3
$ws.Cells.Item(98, 8).Value = 809.6667
$ws.Cells.Item(98, 9).Value = 719.6
$ws.Cells.Item(98, 10).Value = 1260
$ws.Cells.Item(98, 11).Value = 719.6
$ws.Cells.Item(98, 12).Value = 1260
$ws.Cells.Item(98, 13).Value = 778.4
$ws.Cells.Item(98, 14).Value = -4256
$ws.Cells.Item(122, 8).Value = 809.6667
$ws.Cells.Item(122, 9).Value = 719.6
$ws.Cells.Item(122, 10).Value = 1260
$ws.Cells.Item(122, 11).Value = 2158.8
$ws.Cells.Item(122, 12).Value = 3780
$ws.Cells.Item(122, 13).Value = 291.1999999999998
$ws.Cells.Item(122, 14).Value = -8680
$ws.Cells.Item(135, 8).Value = 1918.091
$ws.Cells.Item(135, 9).Value = 642.0526
$ws.Cells.Item(135, 11).Value = 5778.4734
$ws.Cells.Item(135, 13).Value = -3243.4734
$ws.Cells.Item(137, 8).Value = 36874812
$ws.Cells.Item(137, 9).Value = 90910240
$ws.Cells.Item(137, 11).Value = 272730720
$ws.Cells.Item(137, 13).Value = -272728170

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2924.4783
$ws.Cells.Item(2, 9).Value = 1440.9474
$ws.Cells.Item(2, 11).Value = 1440.9474
$ws.Cells.Item(2, 13).Value = -1327.9474
$ws.Cells.Item(74, 8).Value = 3127412.5
$ws.Cells.Item(74, 9).Value = 3788808.8
$ws.Cells.Item(74, 11).Value = 3788808.8
$ws.Cells.Item(74, 13).Value = -3787934.8
$ws.Cells.Item(77, 8).Value = 3127412.5
$ws.Cells.Item(77, 9).Value = 3788808.8
$ws.Cells.Item(77, 11).Value = 18944044
$ws.Cells.Item(77, 13).Value = -18939676
$ws.Cells.Item(97, 8).Value = 738.05884
$ws.Cells.Item(97, 9).Value = 776.1875
$ws.Cells.Item(97, 11).Value = 776.1875
$ws.Cells.Item(97, 13).Value = -280.1875
$ws.Cells.Item(116, 8).Value = 2924.4783
$ws.Cells.Item(116, 9).Value = 1440.9474
$ws.Cells.Item(116, 11).Value = 1440.9474
$ws.Cells.Item(116, 13).Value = 853.0526
$ws.Cells.Item(132, 8).Value = 553531.8
$ws.Cells.Item(132, 9).Value = 678645.5
$ws.Cells.Item(132, 11).Value = 2035936.5
$ws.Cells.Item(132, 13).Value = -2033406.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2924.4783
$ws.Cells.Item(3, 9).Value = 1440.9474
$ws.Cells.Item(3, 11).Value = 1440.9474
$ws.Cells.Item(3, 13).Value = -1326.9474
$ws.Cells.Item(94, 8).Value = 829.34485
$ws.Cells.Item(94, 9).Value = 742.6667
$ws.Cells.Item(94, 10).Value = 971.1818
$ws.Cells.Item(94, 11).Value = 742.6667
$ws.Cells.Item(94, 12).Value = 971.1818
$ws.Cells.Item(94, 13).Value = -291.6667
$ws.Cells.Item(94, 14).Value = -1873.1818
$ws.Cells.Item(134, 8).Value = 628556.5600000001
$ws.Cells.Item(134, 9).Value = 864566.7
$ws.Cells.Item(134, 11).Value = 2593700.1
$ws.Cells.Item(134, 13).Value = -2591165.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 569.375
$ws.Cells.Item(10, 9).Value = 553.5
$ws.Cells.Item(10, 10).Value = 574.6667
$ws.Cells.Item(10, 11).Value = 553.5
$ws.Cells.Item(10, 12).Value = 574.6667
$ws.Cells.Item(10, 13).Value = -414.5
$ws.Cells.Item(10, 14).Value = -852.6667
$ws.Cells.Item(31, 8).Value = 81434.57000000001
$ws.Cells.Item(31, 9).Value = 121317.734
$ws.Cells.Item(31, 10).Value = 23825.555
$ws.Cells.Item(31, 11).Value = 121317.734
$ws.Cells.Item(31, 12).Value = 23825.555
$ws.Cells.Item(31, 13).Value = -121022.734
$ws.Cells.Item(31, 14).Value = -24415.555
$ws.Cells.Item(34, 8).Value = 81434.57000000001
$ws.Cells.Item(34, 9).Value = 121317.734
$ws.Cells.Item(34, 10).Value = 23825.555
$ws.Cells.Item(34, 11).Value = 121317.734
$ws.Cells.Item(34, 12).Value = 23825.555
$ws.Cells.Item(34, 13).Value = -121115.734
$ws.Cells.Item(34, 14).Value = -24229.555
$ws.Cells.Item(88, 8).Value = 45340
$ws.Cells.Item(88, 10).Value = 59950
$ws.Cells.Item(88, 12).Value = 59950
$ws.Cells.Item(88, 14).Value = -60762
$ws.Cells.Item(91, 8).Value = 45340
$ws.Cells.Item(91, 10).Value = 59950
$ws.Cells.Item(91, 12).Value = 59950
$ws.Cells.Item(91, 14).Value = -62758
$ws.Cells.Item(120, 8).Value = 69666.664
$ws.Cells.Item(120, 10).Value = 69666.664
$ws.Cells.Item(120, 12).Value = 69666.664
$ws.Cells.Item(120, 14).Value = -76924.664
$ws.Cells.Item(132, 8).Value = 33039436
$ws.Cells.Item(132, 9).Value = 40017068
$ws.Cells.Item(132, 10).Value = 19620912
$ws.Cells.Item(132, 11).Value = 120051204
$ws.Cells.Item(132, 12).Value = 58862736
$ws.Cells.Item(132, 13).Value = -120048674
$ws.Cells.Item(132, 14).Value = -58867796
$ws.Cells.Item(134, 8).Value = 3589843.8
$ws.Cells.Item(134, 9).Value = 24622.732
$ws.Cells.Item(134, 10).Value = 12502897
$ws.Cells.Item(134, 11).Value = 73868.196
$ws.Cells.Item(134, 12).Value = 37508691
$ws.Cells.Item(134, 13).Value = -71333.196
$ws.Cells.Item(134, 14).Value = -37513761

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1322.8334
$ws.Cells.Item(5, 9).Value = 587.4
$ws.Cells.Item(5, 11).Value = 1762.2
$ws.Cells.Item(5, 13).Value = -1650.2
$ws.Cells.Item(14, 8).Value = 200691.9
$ws.Cells.Item(14, 9).Value = 200691.9
$ws.Cells.Item(14, 11).Value = 602075.7
$ws.Cells.Item(14, 13).Value = -601902.7
$ws.Cells.Item(107, 8).Value = 545
$ws.Cells.Item(107, 10).Value = 545
$ws.Cells.Item(107, 12).Value = 1635
$ws.Cells.Item(107, 14).Value = -5475
$ws.Cells.Item(135, 8).Value = 1322.8334
$ws.Cells.Item(135, 9).Value = 587.4
$ws.Cells.Item(135, 11).Value = 5286.599999999999
$ws.Cells.Item(135, 13).Value = -2751.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7600
$ws.Cells.Item(70, 9).Value = 7800
$ws.Cells.Item(70, 11).Value = 7800
$ws.Cells.Item(70, 13).Value = -7530
$ws.Cells.Item(73, 8).Value = 7600
$ws.Cells.Item(73, 9).Value = 7800
$ws.Cells.Item(73, 11).Value = 7800
$ws.Cells.Item(73, 13).Value = -6864
$ws.Cells.Item(80, 8).Value = 2281.2812
$ws.Cells.Item(80, 9).Value = 2297.0527
$ws.Cells.Item(80, 10).Value = 2258.2307
$ws.Cells.Item(80, 11).Value = 2297.0527
$ws.Cells.Item(80, 12).Value = 2258.2307
$ws.Cells.Item(80, 13).Value = -1299.0527
$ws.Cells.Item(80, 14).Value = -4254.2307
$ws.Cells.Item(83, 8).Value = 2281.2812
$ws.Cells.Item(83, 9).Value = 2297.0527
$ws.Cells.Item(83, 10).Value = 2258.2307
$ws.Cells.Item(83, 11).Value = 11485.2635
$ws.Cells.Item(83, 12).Value = 11291.1535
$ws.Cells.Item(83, 13).Value = -6493.263500000001
$ws.Cells.Item(83, 14).Value = -21275.1535
$ws.Cells.Item(96, 8).Value = 34082.332
$ws.Cells.Item(96, 10).Value = 34082.332
$ws.Cells.Item(96, 12).Value = 34082.332
$ws.Cells.Item(96, 14).Value = -39574.332
$ws.Cells.Item(97, 8).Value = 1444.1852
$ws.Cells.Item(97, 9).Value = 1422.8572
$ws.Cells.Item(97, 11).Value = 1422.8572
$ws.Cells.Item(97, 13).Value = -926.8571999999999
$ws.Cells.Item(122, 8).Value = 27518.441
$ws.Cells.Item(122, 9).Value = 35551.734
$ws.Cells.Item(122, 10).Value = 8980.076999999999
$ws.Cells.Item(122, 11).Value = 106655.202
$ws.Cells.Item(122, 12).Value = 26940.231
$ws.Cells.Item(122, 13).Value = -104205.202
$ws.Cells.Item(122, 14).Value = -31840.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 364.4
$ws.Cells.Item(55, 9).Value = 536
$ws.Cells.Item(55, 10).Value = 250
$ws.Cells.Item(55, 11).Value = 536
$ws.Cells.Item(55, 12).Value = 250
$ws.Cells.Item(55, 13).Value = -363
$ws.Cells.Item(55, 14).Value = -596
$ws.Cells.Item(82, 8).Value = 1016.2353
$ws.Cells.Item(82, 9).Value = 509.33334
$ws.Cells.Item(82, 10).Value = 1292.7273
$ws.Cells.Item(82, 11).Value = 509.33334
$ws.Cells.Item(82, 12).Value = 1292.7273
$ws.Cells.Item(82, 13).Value = -148.33334
$ws.Cells.Item(82, 14).Value = -2014.7273
$ws.Cells.Item(85, 8).Value = 1016.2353
$ws.Cells.Item(85, 9).Value = 509.33334
$ws.Cells.Item(85, 10).Value = 1292.7273
$ws.Cells.Item(85, 11).Value = 509.33334
$ws.Cells.Item(85, 12).Value = 1292.7273
$ws.Cells.Item(85, 13).Value = 738.66666
$ws.Cells.Item(85, 14).Value = -3788.7273

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 82497.5
$ws.Cells.Item(7, 10).Value = 82497.5
$ws.Cells.Item(7, 12).Value = 82497.5
$ws.Cells.Item(7, 14).Value = -82723.5
$ws.Cells.Item(81, 8).Value = 2629.389
$ws.Cells.Item(81, 9).Value = 1384.9166
$ws.Cells.Item(81, 11).Value = 2769.8332
$ws.Cells.Item(81, 13).Value = -1708.8332
$ws.Cells.Item(84, 8).Value = 2629.389
$ws.Cells.Item(84, 9).Value = 1384.9166
$ws.Cells.Item(84, 11).Value = 13849.166
$ws.Cells.Item(84, 13).Value = -8545.166000000001
$ws.Cells.Item(107, 8).Value = 2496.353
$ws.Cells.Item(107, 9).Value = 1249.6923
$ws.Cells.Item(107, 11).Value = 3749.0769
$ws.Cells.Item(107, 13).Value = -1829.0769
$ws.Cells.Item(132, 8).Value = 12946171
$ws.Cells.Item(132, 9).Value = 17446656
$ws.Cells.Item(132, 11).Value = 52339968
$ws.Cells.Item(132, 13).Value = -52337438
$ws.Cells.Item(136, 8).Value = 59497116
$ws.Cells.Item(136, 9).Value = 79287450
$ws.Cells.Item(136, 10).Value = 126115
$ws.Cells.Item(136, 11).Value = 237862350
$ws.Cells.Item(136, 12).Value = 378345
$ws.Cells.Item(136, 13).Value = -237859800
$ws.Cells.Item(136, 14).Value = -383445

